# Iraq League base update (24-02-2024 12:40)
# - Rows that were re-sorted in the source feed end up swapped pairwise
#   (everything except the running "id" in column A moves together).
# - The last match in the sheet (row 211) is replaced by a newer match,
#   and the match that used to be last (row 212) is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Pairs of data rows whose content (columns B..AC) is swapped; the id in
# column A stays where it is.
Swap-Rows 14 15
Swap-Rows 16 17
Swap-Rows 80 81
Swap-Rows 172 173
Swap-Rows 186 187
Swap-Rows 191 192

# Row 211 becomes a new match (id 209 is kept).
$ws.Range("B211").Value2 = 7864432
$ws.Range("E211").Value2 = 45346.54166666666
$ws.Range("F211").Value2 = "Al Shorta SC"
$ws.Range("G211").Value2 = "Al Kahrabaa"
$ws.Range("N211").Value2 = 1.6
$ws.Range("O211").Value2 = 3.3
$ws.Range("P211").Value2 = 5.75
$ws.Range("R211").Value2 = 1.775
$ws.Range("S211").Value2 = 2.025
$ws.Range("U211").Value2 = 1.875
$ws.Range("V211").Value2 = 1.925

# The match that used to be row 212 is removed outright (the sheet shrinks
# from 212 rows to 211 rows).
$ws.Rows.Item(212).Delete()
